# Table writer now honours the incoming column widths instead of always
# dividing the available width evenly across columns. Re-apply the
# (now slightly different, due to rounding) column widths to the single
# table on the slide.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
            $tbl.Columns.Item($c).Width = 198
        }
    }
}
